# Generate Report for Handoff
# Update "Latest Handoff Date(time)" columns to reflect the new handoff pass.
#
# Overview sheet: column D ("Latest Handoff Date") for rows whose status is
# "Ready for handoff" / "Handback transform failed" now shows the new
# handoff timestamp 2016-03-25 09:12:54.
#
# zh-cn sheet: column E ("Latest Handoff Datetime") for the same rows now
# shows 2016-03-25 09:12:44.
#
# de-de sheet: column E ("Latest Handoff Datetime") for the same rows now
# shows 2016-03-25 09:12:54.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$overviewDate = "2016-03-25 09:12:54"
$zhcnDate     = "2016-03-25 09:12:44"
$dedeDate     = "2016-03-25 09:12:54"

# Overview: rows 4, 6, 7, 8, 9, 10 -> column D
$overview.Range("D4").Value  = $overviewDate
$overview.Range("D6").Value  = $overviewDate
$overview.Range("D7").Value  = $overviewDate
$overview.Range("D8").Value  = $overviewDate
$overview.Range("D9").Value  = $overviewDate
$overview.Range("D10").Value = $overviewDate

# zh-cn: rows 4, 6, 7, 8, 9, 10 -> column E
$zhcn.Range("E4").Value  = $zhcnDate
$zhcn.Range("E6").Value  = $zhcnDate
$zhcn.Range("E7").Value  = $zhcnDate
$zhcn.Range("E8").Value  = $zhcnDate
$zhcn.Range("E9").Value  = $zhcnDate
$zhcn.Range("E10").Value = $zhcnDate

# de-de: rows 4, 6, 7, 8, 9, 10 -> column E
$dede.Range("E4").Value  = $dedeDate
$dede.Range("E6").Value  = $dedeDate
$dede.Range("E7").Value  = $dedeDate
$dede.Range("E8").Value  = $dedeDate
$dede.Range("E9").Value  = $dedeDate
$dede.Range("E10").Value = $dedeDate
